# DaySale report: a new sale line (item 6, ZURCAL 40MG) was recorded.
# This pushes the "grand total" row and the footer row down by one row,
# bumps the grand total by the new line's price, and refreshes the
# printed generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Move the footer row (old row 13) down to row 14.
$ws.Range("A13:Q13").Copy($ws.Range("A14:Q14"))
$ws.Rows(14).RowHeight = 16.5
$ws.Range("A13:Q13").UnMerge()
$ws.Range("A13:Q13").Clear()

# 2) Move the grand-total row (old row 12) down to row 13; it gets a
#    slightly taller row in the new layout.
$ws.Range("P12:Q12").Copy($ws.Range("P13:Q13"))
$ws.Rows(13).RowHeight = 24.75

# 3) Build the new sale line in row 12 from the same layout/format as
#    the row above it (row 11), then overwrite with the new item data.
$ws.Range("A12:Q12").UnMerge()
$ws.Range("A12:Q12").Clear()
$ws.Range("A11:Q11").Copy($ws.Range("A12:Q12"))
$ws.Rows(12).RowHeight = 25.5

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "ZURCAL 40MG 14 GASTRO RESISTANT TAB"
$ws.Range("H12").Value = "5:0"

# L12 and P12 carry a numeric display format, but the sheet always
# stores these columns as literal text - coerce text entry explicitly
# so they don't get reinterpreted as numbers.
$fmt = $ws.Range("L12").NumberFormat
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = "1"
$ws.Range("L12").NumberFormat = $fmt

$ws.Range("N12").Value = "96.00"

$fmt = $ws.Range("P12").NumberFormat
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "96.0000"
$ws.Range("P12").NumberFormat = $fmt

$ws.Range("Q12").Value = "1:0"

# 4) Grand total grows by the new line's price (184.66 -> 280.66).
$ws.Range("P13").Value = 280.66000000000003

# 5) Footer timestamp refreshed to the new save time.
$ws.Range("A14").Value = "Saturday, 6 September, 2025 10:05 AM"

Write-Output "done"
